$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '57.377.08'
Set-TextCell 2 5 '  -0.60%  '

# Row 3
Set-TextCell 3 4 '3.101.32'
Set-TextCell 3 5 '  +0.01%  '

# Row 4
Set-TextCell 4 5 '  -0.02%  '

# Row 5
Set-TextCell 5 4 '525.24'
Set-TextCell 5 5 '  +0.11%  '

# Row 6
Set-TextCell 6 4 '136.04'
Set-TextCell 6 5 '  -4.17%  '

# Row 7
Set-TextCell 7 5 '  +0.05%  '

# Row 8
Set-TextCell 8 4 '3.097.83'
Set-TextCell 8 5 '  -0.06%  '

# Row 9
Set-TextCell 9 5 '  +2.29%  '

# Row 10
Set-TextCell 10 4 '7.32'
Set-TextCell 10 5 '  +1.28%  '

# Row 11
Set-TextCell 11 5 '  -1.29%  '

# Row 12
Set-TextCell 12 4 '0.396'
Set-TextCell 12 5 '  +1.75%  '

# Row 13
Set-TextCell 13 4 '3.632.97'
Set-TextCell 13 5 '  -0.13%  '

# Row 14
Set-TextCell 14 5 '  +2.20%  '

# Row 15
Set-TextCell 15 4 '25.26'
Set-TextCell 15 5 '  -2.01%  '

# Row 16
Set-TextCell 16 5 '  -1.11%  '

# Row 17
Set-TextCell 17 4 '57.470.94'
Set-TextCell 17 5 '  -0.62%  '

# Row 18
Set-TextCell 18 4 '3.096.85'
Set-TextCell 18 5 '  -0.16%  '

# Row 19
Set-TextCell 19 4 '5.91'
Set-TextCell 19 5 '  -3.16%  '

# Row 20
Set-TextCell 20 4 '12.35'
Set-TextCell 20 5 '  -3.54%  '

# Row 21
Set-TextCell 21 4 '7.84'
Set-TextCell 21 5 '  -2.50%  '

# Row 22
Set-TextCell 22 4 '347.45'
Set-TextCell 22 5 '  +1.84%  '

# Row 23
Set-TextCell 23 5 '  -0.03%  '

# Row 24
Set-TextCell 24 4 '67.50'
Set-TextCell 24 5 '  +0.96%  '

# Row 25
Set-TextCell 25 5 '  -2.91%  '

# Row 26
Set-TextCell 26 2 'Binance-PegBSC-USD'
Set-TextCell 26 3 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextCell 26 4 '1.02'
Set-TextCell 26 5 '  +1.94%  '

# Row 27
Set-TextCell 27 2 'Kaspa'
Set-TextCell 27 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 27 4 '0.166'
Set-TextCell 27 5 '  -2.13%  '

# Row 28
Set-TextCell 28 4 '0.0₃0890'
Set-TextCell 28 5 '  -3.23%  '

# Row 29
Set-TextCell 29 5 '  -0.04%  '

# Row 30
Set-TextCell 30 4 '7.37'
Set-TextCell 30 5 '  +2.15%  '

# Row 31
Set-TextCell 31 5 '  -0.28%  '

# Row 32
Set-TextCell 32 5 '  -7.52%  '

# Row 33
Set-TextCell 33 4 '20.66'
Set-TextCell 33 5 '  -1.66%  '

# Row 34
Set-TextCell 34 5 '  +6.54%  '

# Row 35
Set-TextCell 35 4 '1.15'
Set-TextCell 35 5 '  -3.93%  '

# Row 36
Set-TextCell 36 4 '158.44'
Set-TextCell 36 5 '  +1.72%  '

# Row 37
Set-TextCell 37 5 '  -1.81%  '

# Row 38
Set-TextCell 38 5 '  -5.17%  '

# Row 39
Set-TextCell 39 5 '  -2.14%  '

# Row 40
Set-TextCell 40 4 '1.61'
Set-TextCell 40 5 '  +6.50%  '

# Row 41
Set-TextCell 41 5 '  -0.98%  '

# Row 42
Set-TextCell 42 5 '  +2.66%  '

# Row 43
Set-TextCell 43 5 '  +1.93%  '

# Row 44
Set-TextCell 44 4 '2.379.83'
Set-TextCell 44 5 '  +3.69%  '

# Row 45
Set-TextCell 45 5 '  -0.80%  '

# Row 46
Set-TextCell 46 5 '  +0.02%  '

# Row 47
Set-TextCell 47 5 '  +1.61%  '

# Row 48
Set-TextCell 48 4 '0.966'
Set-TextCell 48 5 '  -1.76%  '

# Row 49
Set-TextCell 49 5 '  -1.58%  '

# Row 50
Set-TextCell 50 5 '  -4.44%  '

# Row 51
Set-TextCell 51 4 '0.757'
Set-TextCell 51 5 '  +2.67%  '
